# Generate Report for Archive
#
# The localization status for this file moved from "Ready for handoff" to
# "In Translation". That shared string is rendered in four places:
#   - Overview!E2  (zh-cn status column)
#   - Overview!F2  (de-de status column)
#   - zh-cn!C2     (Status column)
#   - de-de!C2     (Status column)
# Updating the text makes the (now shorter) status columns re-fit to their
# new content, so we narrow those columns to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn status in column E, de-de status in column F
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# zh-cn / de-de sheets: Status in column C
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Narrow the status columns now that the text is shorter
$newWidth = 12.5
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
